$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep text formatting so numeric-looking strings are not
# auto-converted to numbers by Excel (matches original inlineStr cell type).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.889.99"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.68%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.886.14"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.41%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -1.20%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.69"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.65%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.03%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4581"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -4.03%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3918"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.16%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "49.16"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -8.94%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08208"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.01%  "
# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.73%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.76"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.14%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.934.73"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.29%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.286"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.46%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.958"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.03%  "
# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.05%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.66"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.82%  "
# Row 18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.81%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06573"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.87%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.32"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -7.23%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.11%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.616"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.61%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.914.06"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.68%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.03"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.47%  "
# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.73%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.127.42"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.90%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.13"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.56%  "
# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.82%  "
# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.42%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.096"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.80%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "123.00"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.95%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09518"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.69%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9533"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.79%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.468"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.20%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.635"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.93%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.434"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.02%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02276"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.41%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.246"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.90%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06084"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.39%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.552"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.84%  "
# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.31%  "
# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.06%  "
# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.62%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1894"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.49%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.300"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.46%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5802"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.42%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.70"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.13%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.981"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.91%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.422"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.40%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06896"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.04%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "110.02"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.24%  "
